$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("NumberError")

# New "Date" column (B) timestamps for rows 2-17, replacing the previous
# batch of run timestamps with a fresh batch from a later test run.
$newDates = @(
    "Mon Oct 09 22:37:12 EDT 2023",
    "Mon Oct 09 22:37:22 EDT 2023",
    "Mon Oct 09 22:37:31 EDT 2023",
    "Mon Oct 09 22:37:41 EDT 2023",
    "Mon Oct 09 22:37:51 EDT 2023",
    "Mon Oct 09 22:38:01 EDT 2023",
    "Mon Oct 09 22:38:11 EDT 2023",
    "Mon Oct 09 22:38:20 EDT 2023",
    "Mon Oct 09 22:38:30 EDT 2023",
    "Mon Oct 09 22:38:39 EDT 2023",
    "Mon Oct 09 22:38:49 EDT 2023",
    "Mon Oct 09 22:38:59 EDT 2023",
    "Mon Oct 09 22:39:08 EDT 2023",
    "Mon Oct 09 22:39:18 EDT 2023",
    "Mon Oct 09 22:39:28 EDT 2023",
    "Mon Oct 09 22:39:37 EDT 2023"
)

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newDates[$i]
}
